$d = $word.ActiveDocument

# =====================================================================
# 1. Summary table at top: "TFS25243 - Add search by Log Name;"
#                       -> "TFS25464 - Update Help url;"
# =====================================================================
$t1 = $d.Tables.Item(1)
$t1.Range.Find.Execute("5243", $true, $false, $false, $false, $false, $true, 1, $false, "5464", 2) | Out-Null
$d.Content.Find.Execute("Add search by Log Name", $true, $false, $false, $false, $false, $true, 1, $false, "Update Help url", 2) | Out-Null

# =====================================================================
# 2. Revision history table: append a new row after the "TFS 25243" row
#    Date: 10/13/2022 | Desc: TFS 25464 - Update Help url | Author: Lili Huang
# =====================================================================
$t = $d.Tables.Item(2)
$newRow = $t.Rows.Add()
$c3 = $newRow.Cells.Item(3)
$c4 = $newRow.Cells.Item(4)
$c3.Merge($c4)

$cell1 = $newRow.Cells.Item(1)
$cell1.Width = 59.3
$cell2 = $newRow.Cells.Item(2)
$cell2.Width = 297.45
$cell3 = $newRow.Cells.Item(3)

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$cell1Xml = $xmlHeader + '<w:p><w:pPr><w:pStyle w:val="TableText"/><w:ind w:left="90"/><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>10/13/2022</w:t></w:r></w:p>' + $xmlFooter
$cell1.Range.InsertXML($cell1Xml)

$cell2Xml = $xmlHeader + '<w:p><w:r><w:t>TFS 25464 – Update Help url</w:t></w:r></w:p>' + $xmlFooter
$cell2.Range.InsertXML($cell2Xml)

$cell3Xml = $xmlHeader + '<w:p><w:pPr><w:pStyle w:val="TableText"/><w:jc w:val="both"/><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Lili Huang</w:t></w:r></w:p>' + $xmlFooter
$cell3.Range.InsertXML($cell3Xml)

Write-Host "New row cell1:" $newRow.Cells.Item(1).Range.Text
Write-Host "New row cell2:" $newRow.Cells.Item(2).Range.Text
Write-Host "New row cell3:" $newRow.Cells.Item(3).Range.Text

# =====================================================================
# 3. Changeset number: 52007 -> 52098
# =====================================================================
$d.Content.Find.Execute("52007", $true, $false, $false, $false, $false, $true, 1, $false, "52098", 2) | Out-Null
